$wb = $excel.ActiveWorkbook

$wsImport = $wb.Worksheets.Item("Import Priorities")
$wsCoupling = $wb.Worksheets.Item("Coupling Parameters")

# --- Update the "type of profit" parameter row (B14) and its description (C14) ---
$wsCoupling.Range("B14").Value = "totalProfits"
$wsCoupling.Range("C14").Value = "totalProfits or irr or none. Total profits are the operational profits excluding the loans"

# Re-style B14: drop the yellow highlight, use a muted monospace "placeholder" look
$wsCoupling.Range("B14").ClearFormats()
$wsCoupling.Range("B14").Font.Name = "JetBrains Mono"
$wsCoupling.Range("B14").Font.Size = 10
$wsCoupling.Range("B14").Font.Color = 8421504
$wsCoupling.Range("B14").VerticalAlignment = -4108

# --- Page setup for the Coupling Parameters sheet ---
$wsCoupling.PageSetup.PaperSize = 9
$wsCoupling.PageSetup.Orientation = 1

# --- Switch the active/selected sheet + selection to Coupling Parameters ---
[void]$wsCoupling.Activate()
[void]$wsCoupling.Range("C13").Select()

# Import Priorities keeps its own (unselected) cursor position
[void]$wsImport.Range("F8").Select()
